$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-04"
$ws.Cells.Item($row, 2).Value = "11:30:58"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "00"

$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 127519
$ws.Cells.Item($row, 6).Value = 143691
$ws.Cells.Item($row, 7).Value = 168186
$ws.Cells.Item($row, 8).Value = 158237
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142077
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192377
$ws.Cells.Item($row, 14).Value = 114749
$ws.Cells.Item($row, 15).Value = 45407
$ws.Cells.Item($row, 16).Value = 28275
$ws.Cells.Item($row, 17).Value = 63193
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47882
$ws.Cells.Item($row, 20).Value = -1
